$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

$ws.Range("H129").Value = 794.5833
$ws.Range("J129").Value = 929.75
$ws.Range("L129").Value = 2789.25
$ws.Range("N129").Value = -12789.25

$ws.Range("H137").Value = 19233102
$ws.Range("I137").Value = 1367.0358
$ws.Range("J137").Value = 41670124
$ws.Range("K137").Value = 4101.107400000001
$ws.Range("L137").Value = 125010372
$ws.Range("M137").Value = -1551.107400000001
$ws.Range("N137").Value = -125015472

$ws.Range("H138").Value = 3133.6711
$ws.Range("I138").Value = 2470.95
$ws.Range("J138").Value = 3870.0278
$ws.Range("K138").Value = 7412.849999999999
$ws.Range("L138").Value = 11610.0834
$ws.Range("M138").Value = -2272.849999999999
$ws.Range("N138").Value = -21890.0834


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14249.833
$ws.Range("I32").Value = 11871.446
$ws.Range("J32").Value = 58250
$ws.Range("K32").Value = 11871.446
$ws.Range("L32").Value = 58250
$ws.Range("M32").Value = -11584.446
$ws.Range("N32").Value = -58824

$ws.Range("H61").Value = 4118962
$ws.Range("I61").Value = 6947699
$ws.Range("J61").Value = 4435.8184
$ws.Range("K61").Value = 6947699
$ws.Range("L61").Value = 4435.8184
$ws.Range("M61").Value = -6947487
$ws.Range("N61").Value = -4859.8184

$ws.Range("H74").Value = 15631027
$ws.Range("I74").Value = 27779078
$ws.Range("J74").Value = 12104.857
$ws.Range("K74").Value = 27779078
$ws.Range("L74").Value = 12104.857
$ws.Range("M74").Value = -27778204
$ws.Range("N74").Value = -13852.857

$ws.Range("H77").Value = 15631027
$ws.Range("I77").Value = 27779078
$ws.Range("J77").Value = 12104.857
$ws.Range("K77").Value = 138895390
$ws.Range("L77").Value = 60524.285
$ws.Range("M77").Value = -138891022
$ws.Range("N77").Value = -69260.285

$ws.Range("H136").Value = 4118962
$ws.Range("I136").Value = 6947699
$ws.Range("J136").Value = 4435.8184
$ws.Range("K136").Value = 20843097
$ws.Range("L136").Value = 13307.4552
$ws.Range("M136").Value = -20840547
$ws.Range("N136").Value = -18407.4552

$ws.Range("H139").Value = 62501.25
$ws.Range("J139").Value = 62501.25
$ws.Range("L139").Value = 62501.25
$ws.Range("N139").Value = -72781.25


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 50860
$ws.Range("J55").Value = 50860
$ws.Range("L55").Value = 50860
$ws.Range("N55").Value = -51406

$ws.Range("H105").Value = 1781.6666
$ws.Range("I105").Value = 1554.2858
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 1554.2858
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = 192.7141999999999
$ws.Range("N105").Value = -5594


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11908651
$ws.Range("I31").Value = 1913.1428
$ws.Range("J31").Value = 23815388
$ws.Range("K31").Value = 1913.1428
$ws.Range("L31").Value = 23815388
$ws.Range("M31").Value = -1618.1428
$ws.Range("N31").Value = -23815978

$ws.Range("H34").Value = 11908651
$ws.Range("I34").Value = 1913.1428
$ws.Range("J34").Value = 23815388
$ws.Range("K34").Value = 1913.1428
$ws.Range("L34").Value = 23815388
$ws.Range("M34").Value = -1711.1428
$ws.Range("N34").Value = -23815792

$ws.Range("H58").Value = 2816.2239
$ws.Range("I58").Value = 1179.1177
$ws.Range("J58").Value = 3372.84
$ws.Range("K58").Value = 1179.1177
$ws.Range("L58").Value = 3372.84
$ws.Range("M58").Value = -976.1177
$ws.Range("N58").Value = -3778.84

$ws.Range("H74").Value = 13299.637
$ws.Range("J74").Value = 15102.889
$ws.Range("L74").Value = 15102.889
$ws.Range("N74").Value = -16850.889

$ws.Range("H77").Value = 13299.637
$ws.Range("J77").Value = 15102.889
$ws.Range("L77").Value = 45308.667
$ws.Range("N77").Value = -54044.667

$ws.Range("H132").Value = 3992.5
$ws.Range("I132").Value = 3768.4546
$ws.Range("K132").Value = 11305.3638
$ws.Range("M132").Value = -8775.363799999999

$ws.Range("H134").Value = 2441.6206
$ws.Range("I134").Value = 2270.85
$ws.Range("J134").Value = 2821.111
$ws.Range("K134").Value = 6812.549999999999
$ws.Range("L134").Value = 8463.332999999999
$ws.Range("M134").Value = -4277.549999999999
$ws.Range("N134").Value = -13533.333

$ws.Range("H136").Value = 2816.2239
$ws.Range("I136").Value = 1179.1177
$ws.Range("J136").Value = 3372.84
$ws.Range("K136").Value = 3537.3531
$ws.Range("L136").Value = 10118.52
$ws.Range("M136").Value = -987.3531000000003
$ws.Range("N136").Value = -15218.52


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 136.34782
$ws.Range("I23").Value = 173.28572
$ws.Range("K23").Value = 519.85716
$ws.Range("M23").Value = -284.85716

$ws.Range("H94").Value = 1000
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null

$ws.Range("H103").Value = 3778808.8
$ws.Range("I103").Value = 11333842
$ws.Range("J103").Value = 1292.3334
$ws.Range("K103").Value = 34001526
$ws.Range("L103").Value = 3877.0002
$ws.Range("M103").Value = -34000647
$ws.Range("N103").Value = -5635.0002

$ws.Range("H109").Value = 2435.2354
$ws.Range("I109").Value = 1350
$ws.Range("J109").Value = 2769.1538
$ws.Range("K109").Value = 4050
$ws.Range("L109").Value = 8307.4614
$ws.Range("M109").Value = -3010
$ws.Range("N109").Value = -10387.4614

$ws.Range("H115").Value = 5930.6924
$ws.Range("I115").Value = 1416.6666
$ws.Range("J115").Value = 9799.857
$ws.Range("K115").Value = 4249.9998
$ws.Range("L115").Value = 29399.571
$ws.Range("M115").Value = -3074.9998
$ws.Range("N115").Value = -31749.571

$ws.Range("H122").Value = 4567081.5
$ws.Range("I122").Value = 13333931
$ws.Range("J122").Value = 1014.4583
$ws.Range("K122").Value = 120005379
$ws.Range("L122").Value = 9130.1247
$ws.Range("M122").Value = -120002929
$ws.Range("N122").Value = -14030.1247

$ws.Range("H138").Value = 3364.0715
$ws.Range("I138").Value = 1656.7858
$ws.Range("J138").Value = 5071.357
$ws.Range("K138").Value = 4970.357400000001
$ws.Range("L138").Value = 15214.071
$ws.Range("M138").Value = 169.6425999999992
$ws.Range("N138").Value = -25494.071


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8907.058999999999
$ws.Range("I80").Value = 2577.1428
$ws.Range("J80").Value = 13338
$ws.Range("K80").Value = 2577.1428
$ws.Range("L80").Value = 13338
$ws.Range("M80").Value = -1579.1428
$ws.Range("N80").Value = -15334

$ws.Range("H83").Value = 8907.058999999999
$ws.Range("I83").Value = 2577.1428
$ws.Range("J83").Value = 13338
$ws.Range("K83").Value = 12885.714
$ws.Range("L83").Value = 66690
$ws.Range("M83").Value = -7893.714
$ws.Range("N83").Value = -76674


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 14998
$ws.Range("I29").Value = 14998
$ws.Range("J29").Value = 14998
$ws.Range("K29").Value = 14998
$ws.Range("L29").Value = 14998
$ws.Range("M29").Value = -14703
$ws.Range("N29").Value = -15588

$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996

$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984

$ws.Range("H132").Value = 7259.84
$ws.Range("I132").Value = 8394.632
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 25183.896
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -22653.896
$ws.Range("N132").Value = -16058.9999

$ws.Range("H136").Value = 1827.2273
$ws.Range("I136").Value = 1124.9375
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 3374.8125
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = -824.8125
$ws.Range("N136").Value = -16200


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12799.2
$ws.Range("I32").Value = 9000
$ws.Range("J32").Value = 15332
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 15332
$ws.Range("M32").Value = -8683
$ws.Range("N32").Value = -15966

$ws.Range("H34").Value = 14998
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 14998
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 14998
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -15404

$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16996

$ws.Range("H81").Value = 4725.7646
$ws.Range("I81").Value = 2299.3333
$ws.Range("J81").Value = 5245.7144
$ws.Range("K81").Value = 4598.6666
$ws.Range("L81").Value = 10491.4288
$ws.Range("M81").Value = -3537.6666
$ws.Range("N81").Value = -12613.4288

$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54984

$ws.Range("H84").Value = 4725.7646
$ws.Range("I84").Value = 2299.3333
$ws.Range("J84").Value = 5245.7144
$ws.Range("K84").Value = 22993.333
$ws.Range("L84").Value = 52457.144
$ws.Range("M84").Value = -17689.333
$ws.Range("N84").Value = -63065.144


Write-Host "Applied all Sheets updates."
